$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card8")

# Row 18 - new service event added to Card8
# Column A: card number (stored as text, like the rest of column A)
$ws.Cells.Item(18, 1).Value = "'8"
$ws.Cells.Item(18, 1).Style = "Normal"

# Columns B-K: left blank, same convention as the rest of the sheet
# (empty text cells rather than truly-empty/untouched cells)
for ($col = 2; $col -le 11; $col++) {
    $ws.Cells.Item(18, $col).Value = "'"
    $ws.Cells.Item(18, $col).Style = "Normal"
}

# Columns L-O: the new event's details
$ws.Cells.Item(18, 12).Value = "20/1/2026"
$ws.Cells.Item(18, 13).Value = "زياره توكيل"
$ws.Cells.Item(18, 14).Value = "تم تغير سوفت كرد لbc"
$ws.Cells.Item(18, 15).Value = "م. احمد علي توكيل"
